$d = $word.ActiveDocument

function Remove-ListParagraph($text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -eq ($text + "`r")) {
            $p.Range.Delete()
            return $true
        }
    }
    return $false
}

# The "FUNCIONALIDADES FALTANTES:" (missing functionalities) bullet list
# currently has three items:
#   1. Modificar Tarjeta de Crédito.
#   2. Desasociar Tarjeta de Crédito.
#   3. Facturar.
#
# Both "Modificar Tarjeta de Crédito" and "Desasociar Tarjeta de Crédito"
# are now implemented, so their list items are removed entirely. The
# remaining "Facturar." item is repurposed to read "Modificar Tarjeta de
# Crédito" (the functionality that is still pending), without a trailing
# period.

Remove-ListParagraph("Modificar Tarjeta de Crédito.")
Remove-ListParagraph("Desasociar Tarjeta de Crédito.")

$d.Content.Find.Execute("Facturar.", $true, $false, $false, $false, $false, $true, 1, $false, "Modificar Tarjeta de Crédito", 2)
